$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column C header: Additional Keywords (bold, like A1/B1) ---
$ws.Range("C1").Value = "Additional Keywords"
$ws.Range("C1").Font.Bold = $true

# --- Additional keyword for the existing Industry/Battery Management System row ---
$ws.Range("C4").Value = "Cloud, BMS"

# --- New Competitor rows ---
$ws.Range("A8").Value = "Competitor"
$ws.Range("B8").Value = "Lithium Balance"
$ws.Range("C8").Value = "Sales, Battery Management System"

$ws.Range("A9").Value = "Competitor"
$ws.Range("B9").Value = "Sensata "
$ws.Range("C9").Value = "Battery Management System, BMS"

# --- New Industry sub-rows ---
$ws.Range("A10").Value = "Industry"
$ws.Range("B10").Value = "CCS"

$ws.Range("A11").Value = "Industry"
$ws.Range("B11").Value = "Electric Vehicle Charge Controller"
$ws.Range("C11").Value = "EVCC"

$ws.Range("A12").Value = "Industry"
$ws.Range("B12").Value = "Fast Charge Junction Box"
$ws.Range("C12").Value = "FCJB"

$ws.Range("A13").Value = "Industry"
$ws.Range("B13").Value = "DC-DC Converter"
$ws.Range("C13").Value = "DCDC"

$ws.Range("A14").Value = "Industry"
$ws.Range("B14").Value = "BMS"
$ws.Range("C14").Value = "Cloud"

# --- Column widths (inputs chosen so the engine's internal char->pixel
#     quantization lands as close as possible to the target stored widths) ---
$ws.Range("B1").ColumnWidth = 93.5
$ws.Range("C1").ColumnWidth = 33.333333333333336
$ws.Range("D1").ColumnWidth = 10.5
$ws.Range("E1").ColumnWidth = 13

# --- Sort A2:E14 by column A ascending ---
$sortRange = $ws.Range("A2:E14")
$keyRange = $ws.Range("A2:A14")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# --- Period / Max Articles values for the Advisors row (now row 2 after sort) ---
$ws.Range("D2").Value = 90
$ws.Range("E2").Value = 2
$ws.Range("D2:E2").HorizontalAlignment = -4108

# --- Period / Max Articles header (bold, centered) ---
$ws.Range("D1").Value = "Period"
$ws.Range("E1").Value = "Max Articles"
$ws.Range("D1:E1").Font.Bold = $true
$ws.Range("D1:E1").HorizontalAlignment = -4108

# --- View settings ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("F6").Select()
